# Update the "Resultaten" section of the measurement report:
#  - the old placeholder sentence is replaced by a reference to the
#    appendix (and a trailing ">" run),
#  - the "_GoBack" bookmark that used to sit at the end of the previous
#    bullet ("Analiseer de resultaten ") is moved in between the new
#    sentence and the trailing ">" run.

$d = $word.ActiveDocument

# --- locate the "Resultaten" body paragraph we need to rewrite -------
$resultsPara = $null   # "Geef de meetresultaten overzichtelijk weer..."

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Geef de meetresultaten overzichtelijk weer*") {
        $resultsPara = $p
        break
    }
}

# --- move the _GoBack bookmark off the "Analiseer de resultaten" bullet
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- replace the body text of the "Resultaten" paragraph -------------
$newText = "Voor de resultaten van deze meting verwijzen we u naar bijlage 1. <rest van de informatie over de verschillende tabellen in de bijlage>"

$bodyRange = $resultsPara.Range
[void]$bodyRange.MoveEnd(1, -1)
$bodyRange.Text = $newText

# --- re-add the bookmark right before the final ">" character --------
$fullRange = $resultsPara.Range
[void]$fullRange.MoveEnd(1, -1)
$bmPos = $fullRange.End - 1

$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
